$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.389.34'
$ws.Cells.Item(2, 5).Value = '  +3.99%  '
$ws.Cells.Item(3, 4).Value = '1.804.24'
$ws.Cells.Item(3, 5).Value = '  +1.53%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.001'
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '316.23'
$ws.Cells.Item(5, 5).Value = '  +0.75%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.001'
$ws.Cells.Item(6, 5).Value = '  +0.07%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5487'
$ws.Cells.Item(7, 5).Value = '  +5.14%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3848'
$ws.Cells.Item(8, 5).Value = '  +6.75%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07580'
$ws.Cells.Item(9, 5).Value = '  +2.73%  '
$ws.Cells.Item(10, 2).Value = 'OKB'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '42.43'
$ws.Cells.Item(10, 5).Value = '  -0.30%  '
$ws.Cells.Item(11, 2).Value = 'Polygon'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.126'
$ws.Cells.Item(11, 5).Value = '  +3.16%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.001'
$ws.Cells.Item(12, 5).Value = '  +0.00%  '
$ws.Cells.Item(13, 5).Value = '  +3.05%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.185'
$ws.Cells.Item(14, 5).Value = '  +1.86%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.388'
$ws.Cells.Item(15, 5).Value = '  +5.79%  '
$ws.Cells.Item(16, 4).Value = '1.806.73'
$ws.Cells.Item(16, 5).Value = '  +1.41%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '91.98'
$ws.Cells.Item(17, 5).Value = '  +4.01%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.00001070'
$ws.Cells.Item(18, 5).Value = '  +2.43%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06445'
$ws.Cells.Item(19, 5).Value = '  +0.48%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.9997'
$ws.Cells.Item(20, 5).Value = '  -0.02%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '17.32'
$ws.Cells.Item(21, 5).Value = '  +3.49%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.972'
$ws.Cells.Item(22, 5).Value = '  +2.07%  '
$ws.Cells.Item(23, 4).Value = '28.409.75'
$ws.Cells.Item(23, 5).Value = '  +3.70%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.44'
$ws.Cells.Item(24, 5).Value = '  +1.09%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.131'
$ws.Cells.Item(25, 5).Value = '  +2.98%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '159.13'
$ws.Cells.Item(26, 5).Value = '  +3.60%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '20.65'
$ws.Cells.Item(27, 5).Value = '  +2.84%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.403'
$ws.Cells.Item(28, 5).Value = '  +2.34%  '
$ws.Cells.Item(29, 4).Value = '2.013.27'
$ws.Cells.Item(29, 5).Value = '  +1.30%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '123.86'
$ws.Cells.Item(30, 5).Value = '  +1.96%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.121'
$ws.Cells.Item(31, 5).Value = '  +5.53%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.1019'
$ws.Cells.Item(32, 5).Value = '  +4.12%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.741'
$ws.Cells.Item(33, 5).Value = '  +3.29%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.684'
$ws.Cells.Item(34, 5).Value = '  +2.14%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.2307'
$ws.Cells.Item(35, 5).Value = '  +13.84%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.06453'
$ws.Cells.Item(36, 5).Value = '  +7.97%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.02325'
$ws.Cells.Item(37, 5).Value = '  +3.97%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '8.847'
$ws.Cells.Item(38, 5).Value = '  +9.29%  '
$ws.Cells.Item(39, 5).Value = '  +6.44%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '11.62'
$ws.Cells.Item(40, 5).Value = '  +3.76%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.6418'
$ws.Cells.Item(41, 5).Value = '  +4.36%  '
$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.161'
$ws.Cells.Item(42, 5).Value = '  +0.98%  '
$ws.Cells.Item(43, 2).Value = 'Frax'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.000'
$ws.Cells.Item(43, 5).Value = '  +0.12%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.385'
$ws.Cells.Item(44, 5).Value = '  -2.94%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '13.58'
$ws.Cells.Item(45, 5).Value = '  +3.02%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.5983'
$ws.Cells.Item(46, 5).Value = '  +3.73%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.685'
$ws.Cells.Item(47, 5).Value = '  +1.51%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '126.92'
$ws.Cells.Item(48, 5).Value = '  +4.50%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.984'
$ws.Cells.Item(49, 5).Value = '  +4.87%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.147'
$ws.Cells.Item(50, 5).Value = '  +3.42%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.06899'
$ws.Cells.Item(51, 5).Value = '  +2.70%  '
